# Updates cryptos list prices / 1h-volume percentages (and re-ranks the
# Dai / RenderToken rows) to match the latest scrape.
#
# Column D ("Price") cells are stored as TEXT in the workbook (e.g. the
# European-style thousands-dot notation "51.515.77", or numbers that would
# otherwise read back as plain floats like "104.02"). Assigning a bare
# numeric-looking string to Range.Value lets Excel auto-coerce it to a
# number, so we prefix with an apostrophe (Excel's literal "treat as text"
# marker) and then reset the cell Style back to Normal so no stray
# quote-prefix formatting lingers on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'51.515.77"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.41%  '
$ws.Range('D3').Value = "'3.102.99"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.88%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').Value = "'385.78"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.48%  '
$ws.Range('D6').Value = "'104.02"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.46%  '
$ws.Range('D7').Value = "'0.539"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -1.30%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('E9').Value = '  -2.09%  '
$ws.Range('E10').Value = '  -0.14%  '
$ws.Range('E11').Value = '  +0.07%  '
$ws.Range('D12').Value = "'0.0858"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.59%  '
$ws.Range('D13').Value = "'3.596.72"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.99%  '
$ws.Range('D14').Value = "'18.54"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.49%  '
$ws.Range('E15').Value = '  +0.59%  '
$ws.Range('D16').Value = "'3.094.41"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.89%  '
$ws.Range('D17').Value = "'0.999"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.46%  '
$ws.Range('D18').Value = "'10.91"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +3.59%  '
$ws.Range('D19').Value = "'51.582.01"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.34%  '
$ws.Range('D20').Value = "'3.28"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +7.16%  '
$ws.Range('D21').Value = "'12.50"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.72%  '
$ws.Range('E22').Value = '  -0.28%  '
$ws.Range('D23').Value = "'70.04"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.02%  '
$ws.Range('D24').Value = "'266.69"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.95%  '
$ws.Range('D25').Value = "'3.17"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.11%  '
$ws.Range('E26').Value = '  -1.28%  '
$ws.Range('D27').Value = "'27.05"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.50%  '
$ws.Range('B28').Value = 'RenderToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D28').Value = "'7.23"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -4.60%  '
$ws.Range('B29').Value = 'Dai'
$ws.Range('C29').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D29').Value = "'1.00"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.02%  '
$ws.Range('E30').Value = '  -4.24%  '
$ws.Range('E31').Value = '  -2.35%  '
$ws.Range('D32').Value = "'10.43"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.79%  '
$ws.Range('D33').Value = "'0.0479"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +5.85%  '
$ws.Range('D34').Value = "'35.48"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +3.28%  '
$ws.Range('D35').Value = "'2.06"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.40%  '
$ws.Range('E36').Value = '  -0.91%  '
$ws.Range('D37').Value = "'0.999"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.08%  '
$ws.Range('D38').Value = "'3.37"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.43%  '
$ws.Range('D39').Value = "'0.294"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.96%  '
$ws.Range('E40').Value = '  -0.45%  '
$ws.Range('D41').Value = "'128.96"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.00%  '
$ws.Range('E42').Value = '  -0.47%  '
$ws.Range('D43').Value = "'16.58"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.28%  '
$ws.Range('E44').Value = '  -3.54%  '
$ws.Range('D45').Value = "'3.79"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.16%  '
$ws.Range('D46').Value = "'22.12"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.85%  '
$ws.Range('D47').Value = "'2.53"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +5.91%  '
$ws.Range('E48').Value = '  -2.62%  '
$ws.Range('D49').Value = "'2.075.64"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.77%  '
$ws.Range('D50').Value = "'0.932"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +18.81%  '
$ws.Range('D51').Value = "'0.0329"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.14%  '
